# Sprint Backlog update:
#   - Row 10 previously tracked a single combined task
#     "Identifie 3 design paterns - Martim Costa" in column A.
#   - Split it to match the pattern used in row 4: column A now just
#     holds the assignee's name ("Martim Costa") and column D holds the
#     task description ("Identifie 3 design paterns").
#   - Update the active selection to D10 (previously D20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Martim Costa"
$ws.Range("D10").Value = "Identifie 3 design paterns"

$ws.Range("D10").Select()
